# daily auto push: 2026-02-22 18:52 UTC
# Insert 3 new data rows (before the existing 2026/12/29 block) for the
# latest scraped timestamps: 2026/02/22 19:00, 2026/02/22 22:00, and
# 2026/02/23 02:00. Everything from the old row 851 onward shifts down
# by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 851/852/853, pushing the rest of the table down.
$ws.Range("A851:A853").EntireRow.Insert()

# New row 851: 2026/02/22, 日, 19, 201
$ws.Range("A851").NumberFormat = "@"
$ws.Range("A851").Value = "2026/02/22"
$ws.Range("A851").Style = "Normal"
$ws.Range("B851").Value = "日"
$ws.Range("C851").Value = 19
$ws.Range("D851").Value = 201

# New row 852: 2026/02/22, 日, 22, 201
$ws.Range("A852").NumberFormat = "@"
$ws.Range("A852").Value = "2026/02/22"
$ws.Range("A852").Style = "Normal"
$ws.Range("B852").Value = "日"
$ws.Range("C852").Value = 22
$ws.Range("D852").Value = 201

# New row 853: 2026/02/23, 月, 2, 201
$ws.Range("A853").NumberFormat = "@"
$ws.Range("A853").Value = "2026/02/23"
$ws.Range("A853").Style = "Normal"
$ws.Range("B853").Value = "月"
$ws.Range("C853").Value = 2
$ws.Range("D853").Value = 201
